## ----------------------------------------------------------------------
## "added sts slowing reason(not complete solution) and started adding
##  ajax notes"
##
## 1. Populate a "many-to-one / one-to-many" JPA annotation cheat sheet on
##    the "notes|annotation usage" sheet, plus a stray note on photoshop.
## 2. Add a chunk of notes about why STS (Spring Tool Suite) runs slowly,
##    with the eclipse.ini before/after tuning values, to "Useful Notes".
## 3. Add two new sheets: "Query" (placeholder, still empty) and "Ajax"
##    (notes about fetch()/promises just getting started).
## ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "notes|annotation usage" sheet -- JPA relationship annotations
# ---------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("notes|annotation usage")

$wsNotes.Range("N11").Value = "photoshop --> zeplin for css"

$wsNotes.Range("A44:P44").Copy()
$wsNotes.Range("A49:P49").PasteSpecial(-4122) # xlPasteFormats
$wsNotes.Range("A49").Value = "다대일"
$wsNotes.Range("H49").Value = "일대다"

$wsNotes.Range("A50").Value = "@ManyToOne"
$wsNotes.Range("H50").Value = "@OneToMany(mappedBy=""uservo"", fetch=FetchType.EAGER, cascade=CascadeType.ALL)"

$wsNotes.Range("A51").Value = "@JoinColumn(name=""user_id"", nullable=false)"
$wsNotes.Range("H51").Value = "private List<ReviewRegistrationvo> reviewRegistrationList = new ArrayList<ReviewRegistrationvo>(); "

$wsNotes.Range("A52").Value = "private Uservo uservo;"

$wsNotes.Range("L13").Select()

# ---------------------------------------------------------------------
# 2. "Useful Notes" sheet -- why is STS slow, eclipse.ini tuning
# ---------------------------------------------------------------------
$wsUseful = $wb.Worksheets.Item("Useful Notes")

$wsUseful.Range("A30").Value = "sts running realllly slow, what to do? "

$wsUseful.Range("A45").Value = "-XX:+UseG1GC"
$wsUseful.Range("A46").Value = "-XX:+UseStringDeduplication"
$wsUseful.Range("A47").Style = $wsUseful.Range("A45").Style
$wsUseful.Range("A47").Value = "=--add-modules=ALL-SYSTEM"
$wsUseful.Range("A48").Value = "-javaagent:C:\Users\Yoon Taewon\Downloads\spring-tool-suite-4-4.3.1.RELEASE-e4.12.0-win32.win32.x86_64\sts-4.3.1.RELEASE\lombok.jar"

$wsUseful.Range("A31").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A31").Value = "=-startup"
$wsUseful.Range("A32").Value = "plugins/org.eclipse.equinox.launcher_1.5.400.v20190515-0925.jar"
$wsUseful.Range("A33").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A33").Value = "=-vm"
$wsUseful.Range("A34").Value = "C:/Program Files/Java/jdk1.8.0_221/bin/javaw.exe"
$wsUseful.Range("A35").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A35").Value = "=--launcher.library"
$wsUseful.Range("A36").Value = "plugins/org.eclipse.equinox.launcher.win32.win32.x86_64_1.1.1000.v20190125-2016"
$wsUseful.Range("A37").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A37").Value = "=-product"
$wsUseful.Range("A38").Value = "org.springframework.boot.ide.branding.sts4"
$wsUseful.Range("A39").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A39").Value = "=--launcher.defaultAction"
$wsUseful.Range("A40").Value = "openFile"
$wsUseful.Range("A41").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A41").Value = "=-vmargs"
$wsUseful.Range("A42").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A42").Value = "=-Dosgi.requiredJavaVersion=1.8"
$wsUseful.Range("A43").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A43").Value = "=-Xms512m"
$wsUseful.Range("A44").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A44").Value = "=-Xmx1024m"

$wsUseful.Range("A50").Value = "original file above"
$wsUseful.Range("A51").Value = "now changed part below"

$wsUseful.Range("A53").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A53").Value = "=-Xms512m"
$wsUseful.Range("C53").Value = "xms랑 xmx를 똑같으 크기로 설정, 최소 heap메모리랑 최대heap메모리 변경이 없어서 속도 향상"
$wsUseful.Range("A54").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A54").Value = "=-Xmx512m"

$wsUseful.Range("A56").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A56").Value = "=-Xverify:none"
$wsUseful.Range("A57").Value = "-XX:+AggressiveOpts"
$wsUseful.Range("A58").Value = "-XX:-UseConcMarkSweepGC"
$wsUseful.Range("A59").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A59").Value = "=-Xmn256m"
$wsUseful.Range("A60").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A60").Value = "=-Xms1024m"
$wsUseful.Range("A61").Style = $wsUseful.Range("A47").Style
$wsUseful.Range("A61").Value = "=-Xmx1024m"
$wsUseful.Range("A62").Value = "-XX:PermSize=128m"
$wsUseful.Range("A63").Value = "-XX:MaxPermSize=128m"
$wsUseful.Range("A64").Value = "-XX:NewSize=128m"
$wsUseful.Range("A65").Value = "-XX:MaxNewSize=128m"

$wsUseful.Range("A67").Value = "prefernce--general--show heap status"
$wsUseful.Range("A68").Value = "turn off unused project"

$wsUseful.Range("A70:H70").Style = $wsUseful.Range("A1:H1").Style

$wsUseful.Range("E64").Select()

# ---------------------------------------------------------------------
# 3. New sheets: "Query" (empty placeholder) and "Ajax" (fetch notes)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsQuery = $wb.Worksheets.Add($null, $lastSheet)
$wsQuery.Name = "Query"

$wsAjax = $wb.Worksheets.Add($null, $wsQuery)
$wsAjax.Name = "Ajax"

$wsAjax.Range("A4").Value = "<a onclick="""
$wsAjax.Range("A5").Value = "document.querySelector('article').innerHTML = '<h2>html</h2>html is cool';>"

$wsAjax.Range("A7").Value = "fetch('/').then(function(response){ "
$wsAjax.Range("A8").Value = "response.text().then(function(text){"
$wsAjax.Range("A9").Value = "document.querySelector('article').innerHTML= text;"
$wsAjax.Range("A10").Value = "}) "
$wsAjax.Range("A11").Value = "}) "

$wsAjax.Range("A13").Value = "function callbackme(){"
$wsAjax.Range("E13").Value = "callbackme = function(){"
$wsAjax.Range("B14").Value = "console.log('end');"
$wsAjax.Range("D14").Style = $wsNotes.Range("L16").Style
$wsAjax.Range("D14").Value = "===="
$wsAjax.Range("F14").Value = "console.log('end');"
$wsAjax.Range("A15").Value = "}"
$wsAjax.Range("E15").Value = "}"

$wsAjax.Range("A16").Value = "fetch('/').then(callbackme);"
$wsAjax.Range("D16").Style = $wsNotes.Range("L16").Style
$wsAjax.Range("D16").Value = "== fetch를 한후에 실행을 끝날때가지 기다리고 callbackme란 function를 실행시킨다"

$wsAjax.Range("D17").Select()
